$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B15").Value = 2.2
$ws.Range("C15").Value = 2.3
$ws.Range("F15").Value = 0.6
$ws.Range("C16").Value = 2

$ws.Range("C17").Select()
